$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 385
$ws.Range("I6").Value = 385
$ws.Range("K6").Value = 1155
$ws.Range("M6").Value = -1043

$ws.Range("H48").Value = 2840
$ws.Range("J48").Value = 3510
$ws.Range("L48").Value = 10530
$ws.Range("N48").Value = -11114

$ws.Range("H56").Value = 2840
$ws.Range("J56").Value = 3510
$ws.Range("L56").Value = 10530
$ws.Range("N56").Value = -11598

$ws.Range("H138").Value = 4146.256
$ws.Range("J138").Value = 4997.697
$ws.Range("L138").Value = 14993.091
$ws.Range("N138").Value = -25273.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1309087.6
$ws.Range("I32").Value = 550.6866
$ws.Range("K32").Value = 550.6866
$ws.Range("M32").Value = -263.6866

$ws.Range("H110").Value = 4106.9062
$ws.Range("I110").Value = 1923.85
$ws.Range("K110").Value = 1923.85
$ws.Range("M110").Value = 121.1500000000001

$ws.Range("H132").Value = 712617.9
$ws.Range("I132").Value = 979661.9
$ws.Range("J132").Value = 75820.69500000001
$ws.Range("K132").Value = 2938985.7
$ws.Range("L132").Value = 227462.085
$ws.Range("M132").Value = -2936455.7
$ws.Range("N132").Value = -232522.085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1199.6666
$ws.Range("J5").Value = 1499.5
$ws.Range("L5").Value = 1499.5
$ws.Range("N5").Value = -1725.5

$ws.Range("H105").Value = 28574980
$ws.Range("I105").Value = 111116590
$ws.Range("K105").Value = 111116590
$ws.Range("M105").Value = -111114843

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 41671090
$ws.Range("I16").Value = 83336450
$ws.Range("J16").Value = 5731.1665
$ws.Range("K16").Value = 83336450
$ws.Range("L16").Value = 5731.1665
$ws.Range("M16").Value = -83336163
$ws.Range("N16").Value = -6305.1665

$ws.Range("H19").Value = 381.66666
$ws.Range("I19").Value = 250
$ws.Range("J19").Value = 447.5
$ws.Range("K19").Value = 250
$ws.Range("L19").Value = 447.5
$ws.Range("M19").Value = -80
$ws.Range("N19").Value = -787.5

$ws.Range("H22").Value = 1716609.9
$ws.Range("I22").Value = 2694537.2
$ws.Range("J22").Value = 5237
$ws.Range("K22").Value = 2694537.2
$ws.Range("L22").Value = 5237
$ws.Range("M22").Value = -2694187.2
$ws.Range("N22").Value = -5937

$ws.Range("H23").Value = 19997.5
$ws.Range("J23").Value = 19997.5
$ws.Range("L23").Value = 19997.5
$ws.Range("N23").Value = -20477.5

$ws.Range("H24").Value = 381.66666
$ws.Range("I24").Value = 250
$ws.Range("J24").Value = 447.5
$ws.Range("K24").Value = 250
$ws.Range("L24").Value = 447.5
$ws.Range("M24").Value = -80
$ws.Range("N24").Value = -787.5

$ws.Range("H27").Value = 19997.5
$ws.Range("J27").Value = 19997.5
$ws.Range("L27").Value = 19997.5
$ws.Range("N27").Value = -20381.5

$ws.Range("H105").Value = 52637596
$ws.Range("I105").Value = 100005830
$ws.Range("K105").Value = 100005830
$ws.Range("M105").Value = -100004083

$ws.Range("H113").Value = 41671090
$ws.Range("I113").Value = 83336450
$ws.Range("J113").Value = 5731.1665
$ws.Range("K113").Value = 83336450
$ws.Range("L113").Value = 5731.1665
$ws.Range("M113").Value = -83334280
$ws.Range("N113").Value = -10071.1665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1688.4
$ws.Range("I7").Value = 2152.7273
$ws.Range("J7").Value = 411.5
$ws.Range("K7").Value = 6458.1819
$ws.Range("L7").Value = 1234.5
$ws.Range("M7").Value = -6346.1819
$ws.Range("N7").Value = -1458.5

$ws.Range("H23").Value = 207
$ws.Range("J23").Value = 406.33334
$ws.Range("L23").Value = 1219.00002
$ws.Range("N23").Value = -1689.00002

$ws.Range("H24").Value = 3399.6667
$ws.Range("J24").Value = 4999.5
$ws.Range("L24").Value = 14998.5
$ws.Range("N24").Value = -15458.5

$ws.Range("H68").Value = 5000000
$ws.Range("J68").Value = 5000000
$ws.Range("L68").Value = 15000000
$ws.Range("N68").Value = -15001622

$ws.Range("H71").Value = 5000000
$ws.Range("J71").Value = 5000000
$ws.Range("L71").Value = 45000000
$ws.Range("N71").Value = -45008112

$ws.Range("H92").Value = 441.32
$ws.Range("J92").Value = 224.66667
$ws.Range("L92").Value = 674.00001
$ws.Range("N92").Value = -3170.00001

$ws.Range("H136").Value = 23812226
$ws.Range("I136").Value = 12822984
$ws.Range("J136").Value = 41669744
$ws.Range("K136").Value = 38468952
$ws.Range("L136").Value = 125009232
$ws.Range("M136").Value = -38463852
$ws.Range("N136").Value = -125019432

$ws.Range("H137").Value = 1796.6666
$ws.Range("I137").Value = 945
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 2835
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = 2265
$ws.Range("N137").Value = -20700

$ws.Range("H138").Value = 7412388.5
$ws.Range("I138").Value = 8613.333000000001
$ws.Range("J138").Value = 7941229.5
$ws.Range("K138").Value = 25839.999
$ws.Range("L138").Value = 23823688.5
$ws.Range("M138").Value = -20699.999
$ws.Range("N138").Value = -23833968.5

$ws.Range("H139").Value = 3067.25
$ws.Range("I139").Value = 2098
$ws.Range("K139").Value = 6294
$ws.Range("M139").Value = -1154

$ws.Range("H140").Value = 55557890
$ws.Range("I140").Value = 75758936
$ws.Range("K140").Value = 227276808
$ws.Range("M140").Value = -227271628

$ws.Range("H141").Value = 5387.9
$ws.Range("I141").Value = 5387.9
$ws.Range("K141").Value = 16163.7
$ws.Range("M141").Value = -10983.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 25005300
$ws.Range("J8").Value = 10599.5
$ws.Range("L8").Value = 10599.5
$ws.Range("N8").Value = -10879.5

$ws.Range("H132").Value = 4935.8413
$ws.Range("J132").Value = 24999
$ws.Range("L132").Value = 74997
$ws.Range("N132").Value = -80057
